# CORRIDAS TK MULTIMARCAS.xlsx - add the missing "BARROS" entry (15) on
# "MES 01" for the 2023-10-17 day group (row 59, which was a placeholder
# 0/blank row) and move the view/selection to where that edit was made.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MES 01")

# Row 59 was an empty placeholder (A59=0, B59 blank) under the 45216
# (2023-10-17) date-header row. Fill in the real corrida entry.
$ws.Range("A59").Value = 15
$ws.Range("B59").Value = "BARROS"

# Bring the view to where the new data now lives.
$ws.Activate()
$ws.Range("K50").Select()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 1
